# "have advance use the same codes"
# The "Advance" rows (B34:B49) had a VAR1 column (D) that used placeholder
# codes "AX", "AOE" and "A" for advances caused by an error / result-of-play.
# Update those cells so the Advance rows re-use the same short codes as the
# other columns instead of the separate placeholder strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unknown (row 34): "AX" -> "X"
$ws.Range("D34").Value = "X"

# Error-family advances: "AOE" -> "E" (matches the plain Error code)
$ws.Range("D37").Value = "E"   # Error
$ws.Range("D38").Value = "E"   # Error after Passed Ball
$ws.Range("D39").Value = "E"   # Error after Pickoff
$ws.Range("D40").Value = "E"   # Error after Steal
$ws.Range("D43").Value = "E"   # Passed Ball
$ws.Range("D46").Value = "E"   # Wild Pitch
$ws.Range("D47").Value = "E"   # Error after Wild Pitch

# "Result of a play" advances: "A" -> the relevant code
$ws.Range("D41").Value = "FC"  # Fielder's Choice
$ws.Range("D42").Value = "P"   # Result of a Play
$ws.Range("D48").Value = "FC"  # Safe on Steal, by Fielder's Choice
$ws.Range("D49").Value = "FC"  # Defensive Indifference

# Reflect the view state left behind after scrolling/editing that part of
# the sheet (the frozen-pane selection ends on H39).
$ws.Range("H39").Select()
